$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update index2Sequence (column K, rows 2-37) to the new shared value "E7420".
#    This collapses the per-row E776x/E777x values (and the original E7760 block)
#    down to a single repeated value, which also causes the now-unused shared
#    strings to drop out of the table once nothing references them anymore.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 11).Value = "E7420"
}

# 2. Give column K its own distinct font (Arial 11, black) instead of inheriting
#    the mixed original formatting. Style K2 directly, then propagate that exact
#    format to the rest of the column so every row collapses onto one single
#    style entry instead of several derived ones.
$kFirst = $ws.Range("K2")
$kFirst.Font.Name = "Arial"
$kFirst.Font.Size = 11
$kFirst.Font.Color = 0
$kFirst.Copy()
$null = $ws.Range("K3:K37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Column L: replace the static boolean constant with a live formula that
#    evaluates to FALSE(), for every data row.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}

# 4. Update the sheet's saved selection to K2:K37 (previously L2:L37).
$null = $ws.Range("K2:K37").Select()
